# Update the "Corr/Total" marks on the marksheet (quiz sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row -> Right column (B11): 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row -> Right column (B12): 33 -> 55
$ws.Range("B12").Value = 55

# "Total" row -> Max column (E12): "31/84" -> "55/140"
$ws.Range("E12").Value = "55/140"
